$wb = $excel.ActiveWorkbook

# Overview sheet - G2: Latest HO Xliff Generate Date (shared with de-de Correspond Handoff Datetime)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 05:04:44"

# zh-cn sheet - H2: Correspond Handoff Datetime, K2: Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 05:04:32"
$wsZhCn.Range("K2").Value = "2016-09-06 05:05:31"

# de-de sheet - H2: Correspond Handoff Datetime (shares same value as Overview G2), K2: Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-06 05:04:44"
$wsDeDe.Range("K2").Value = "2016-09-06 05:05:51"
